# Add a new "olac:isDeIdentified" rdf:Property term to the
# "@type=rdfProperty" sheet. The term is inserted alphabetically between
# "olac:isAnnotationOf" and "olac:licensedOrganization", i.e. as the new
# row 19 (pushing the existing row 19 "olac:licensedOrganization" and all
# following rows down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("@type=rdfProperty")

# Insert a blank row before the current row 19.
$ws.Rows.Item(19).Insert()

# @id
$ws.Range("A19").Value = "olac:isDeIdentified"
# @type
$ws.Range("B19").Value = "rdf:Property"
# domainIncludes
$ws.Range("D19").Value = '[{"@id":"schema:CreativeWork"}, {"@id":"schema:Person"}, olac:PersonSnapshot]'
# name
$ws.Range("E19").Value = "isDeIdentified"
# rdfs:label
$ws.Range("F19").Value = "isDeIdentified"
# rangeIncludes
$ws.Range("G19").Value = '{"@id":"schema:Boolean"}'
# rdfs:comment
$ws.Range("I19").Value = "This data in this item has had identifying information removed, or in the case of a person the name is an alias"
